$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''326.62'
$ws.Range("E2").Value = '''0.03%'
$ws.Range("G2").Value = '''4'
$ws.Range("D3").Value = '''44.03'
$ws.Range("E3").Value = '''-2.29%'
$ws.Range("G3").Value = '''4'
$ws.Range("D4").Value = '''5.492'
$ws.Range("E4").Value = '''-1.39%'
$ws.Range("G4").Value = '''4'
$ws.Range("D5").Value = '''0.08023'
$ws.Range("E5").Value = '''-0.75%'
$ws.Range("G5").Value = '''4'
$ws.Range("D6").Value = '''1.975'
$ws.Range("E6").Value = '''4.08%'
$ws.Range("G6").Value = '''4'
$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D7").Value = '''2.567'
$ws.Range("E7").Value = '''-5.14%'
$ws.Range("G7").Value = '''4'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9510'
$ws.Range("E8").Value = '''0.35%'
$ws.Range("G8").Value = '''4'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1132'
$ws.Range("E9").Value = '''-2.29%'
$ws.Range("G9").Value = '''4'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1842'
$ws.Range("E10").Value = '''-2.81%'
$ws.Range("G10").Value = '''4'
$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").Value = '''11.36'
$ws.Range("E11").Value = '''33.74%'
$ws.Range("G11").Value = '''4'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09838'
$ws.Range("E12").Value = '''-3.31%'
$ws.Range("G12").Value = '''4'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04619'
$ws.Range("E13").Value = '''10.45%'
$ws.Range("G13").Value = '''4'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1067'
$ws.Range("E14").Value = '''0.29%'
$ws.Range("G14").Value = '''4'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001270'
$ws.Range("E15").Value = '''-1.58%'
$ws.Range("G15").Value = '''4'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04089'
$ws.Range("E16").Value = '''-4.36%'
$ws.Range("G16").Value = '''4'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.005713'
$ws.Range("E17").Value = '''-3.95%'
$ws.Range("G17").Value = '''4'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.365'
$ws.Range("E18").Value = '''-6.87%'
$ws.Range("G18").Value = '''4'
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").Value = '''4.291'
$ws.Range("E19").Value = '''-0.96%'
$ws.Range("G19").Value = '''4'
$ws.Range("D20").Value = '''0.3479'
$ws.Range("E20").Value = '''-0.18%'
$ws.Range("G20").Value = '''4'
$ws.Range("D21").Value = '''0.1409'
$ws.Range("E21").Value = '''2.06%'
$ws.Range("G21").Value = '''4'
$ws.Range("D22").Value = '''0.2545'
$ws.Range("E22").Value = '''-4.46%'
$ws.Range("G22").Value = '''4'
$ws.Range("D23").Value = '''0.001245'
$ws.Range("E23").Value = '''0.72%'
$ws.Range("G23").Value = '''4'
$ws.Range("E24").Value = '''-7.16%'
$ws.Range("G24").Value = '''4'
$ws.Range("D25").Value = '''0.0001190'
$ws.Range("E25").Value = '''-3.46%'
$ws.Range("G25").Value = '''4'
$ws.Range("D26").Value = '''0.0003741'
$ws.Range("E26").Value = '''-6.45%'
$ws.Range("G26").Value = '''4'
$ws.Range("G27").Value = '''4'
$ws.Range("G28").Value = '''4'
$ws.Range("G29").Value = '''4'
$ws.Range("G30").Value = '''4'
$ws.Range("G31").Value = '''4'
$ws.Range("G32").Value = '''4'
$ws.Range("G33").Value = '''4'
$ws.Range("G34").Value = '''4'
$ws.Range("G35").Value = '''4'
$ws.Range("G36").Value = '''4'
$ws.Range("G37").Value = '''4'
$ws.Range("E38").Value = '''-3.90%'
$ws.Range("G38").Value = '''4'
$ws.Range("D39").Value = '''0.05556'
$ws.Range("E39").Value = '''-0.01%'
$ws.Range("G39").Value = '''4'
$ws.Range("D40").Value = '''0.007544'
$ws.Range("E40").Value = '''-1.47%'
$ws.Range("G40").Value = '''4'
$ws.Range("D41").Value = '''0.1395'
$ws.Range("E41").Value = '''0.09%'
$ws.Range("G41").Value = '''4'
$ws.Range("D42").Value = '''0.007635'
$ws.Range("E42").Value = '''-32.64%'
$ws.Range("G42").Value = '''4'
$ws.Range("D43").Value = '''0.002015'
$ws.Range("E43").Value = '''-2.12%'
$ws.Range("G43").Value = '''4'
$ws.Range("D44").Value = '''0.008494'
$ws.Range("E44").Value = '''-2.05%'
$ws.Range("G44").Value = '''4'
$ws.Range("D45").Value = '''0.00007110'
$ws.Range("E45").Value = '''-0.01%'
$ws.Range("G45").Value = '''4'
$ws.Range("E46").Value = '''-0.33%'
$ws.Range("G46").Value = '''4'
$ws.Range("E47").Value = '''54.99%'
$ws.Range("G47").Value = '''4'
$ws.Range("D48").Value = '''0.003431'
$ws.Range("E48").Value = '''-0.25%'
$ws.Range("G48").Value = '''4'
$ws.Range("E49").Value = '''-0.33%'
$ws.Range("G49").Value = '''4'
$ws.Range("D50").Value = '''0.0001999'
$ws.Range("E50").Value = '''-0.33%'
$ws.Range("G50").Value = '''4'
$ws.Range("G51").Value = '''4'
